# The deck ships with two DrawingML themes:
#   ppt/theme/theme1.xml -> bound to the (only) slide master, currently the
#                            "Integral" design
#   ppt/theme/theme2.xml -> bound to the notes master, currently the
#                            default "Office Theme"
#
# The authored edit swaps the two themes' content: the slide master/theme1
# ends up with the "Office Theme" palette (and theme2 ends up with the
# "Integral" palette). The font scheme and format scheme are identical
# between the two themes already, so only the 12-slot colour scheme (and
# its/​the theme's display name) actually differ.
#
# Apply the new palette through the Theme Colors object model exposed on a
# Slide (it writes straight into the shared theme's <a:clrScheme>, i.e.
# ppt/theme/theme1.xml, since every slide shares the single slide master).

function Set-ThemeColor {
    param(
        $ThemeColorScheme,
        [int]$Index,
        [int]$R,
        [int]$G,
        [int]$B
    )
    # PowerPoint's ColorFormat.RGB uses the OLE COM RGB() packing
    # (0x00BBGGRR), not the 0xRRGGBB order used in the OOXML srgbClr val.
    $bgr = ($B * 65536) + ($G * 256) + $R
    $ThemeColorScheme.Colors($Index).RGB = $bgr
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Office Theme colour scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
Set-ThemeColor $tcs 1  0x00 0x00 0x00   # dk1
Set-ThemeColor $tcs 2  0xFF 0xFF 0xFF   # lt1
Set-ThemeColor $tcs 3  0x44 0x54 0x6A   # dk2
Set-ThemeColor $tcs 4  0xE7 0xE6 0xE6   # lt2
Set-ThemeColor $tcs 5  0x5B 0x9B 0xD5   # accent1
Set-ThemeColor $tcs 6  0xED 0x7D 0x31   # accent2
Set-ThemeColor $tcs 7  0xA5 0xA5 0xA5   # accent3
Set-ThemeColor $tcs 8  0xFF 0xC0 0x00   # accent4
Set-ThemeColor $tcs 9  0x44 0x72 0xC4   # accent5
Set-ThemeColor $tcs 10 0x70 0xAD 0x47   # accent6
Set-ThemeColor $tcs 11 0x05 0x63 0xC1   # hlink
Set-ThemeColor $tcs 12 0x95 0x4F 0x72   # folHlink
